$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B -> old B becomes C, old C becomes D
$ws.Columns.Item(2).Insert()

# Header row
$ws.Range("B1").Value = "Variável"

$varText = "Trabalho como origem na renda (%)"

$dates = @(
    "31/12/2012",
    "31/12/2013",
    "31/12/2014",
    "31/12/2015",
    "31/12/2016",
    "31/12/2017",
    "31/12/2018",
    "31/12/2019"
)

for ($block = 0; $block -lt 3; $block++) {
    for ($i = 0; $i -lt 8; $i++) {
        $row = 2 + ($block * 8) + $i
        $ws.Range("B$row").Value = $varText
        $ws.Range("C$row").Value = $dates[$i]
    }
}
